$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '63.098.28'
$ws.Range("E2").Value = '  +1.27%  '

# Row 3
$ws.Range("D3").Value = '2.473.50'
$ws.Range("E3").Value = '  +1.75%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.08%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '573.21'
$ws.Range("E5").Value = '  +1.67%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '147.23'
$ws.Range("E6").Value = '  +1.71%  '

# Row 7
$ws.Range("E7").Value = '  -0.03%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.541'
$ws.Range("E8").Value = '  +1.78%  '

# Row 9
$ws.Range("D9").Value = '2.473.05'
$ws.Range("E9").Value = '  +1.79%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.112'
$ws.Range("E10").Value = '  +1.41%  '

# Row 11
$ws.Range("E11").Value = '  +0.88%  '

# Row 12
$ws.Range("B12").Value = 'Cardano'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.357'
$ws.Range("E12").Value = '  +1.91%  '

# Row 13
$ws.Range("B13").Value = 'Toncoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.25'
$ws.Range("E13").Value = '  -0.08%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.17'
$ws.Range("E14").Value = '  +1.78%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000179'
$ws.Range("E15").Value = '  +2.59%  '

# Row 16
$ws.Range("D16").Value = '2.910.29'
$ws.Range("E16").Value = '  +2.34%  '

# Row 17
$ws.Range("D17").Value = '62.869.09'
$ws.Range("E17").Value = '  +0.97%  '

# Row 18
$ws.Range("D18").Value = '2.481.23'
$ws.Range("E18").Value = '  +2.16%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.45'
$ws.Range("E19").Value = '  +1.63%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.33'
$ws.Range("E20").Value = '  +6.90%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '324.01'
$ws.Range("E21").Value = '  -0.48%  '

# Row 22
$ws.Range("E22").Value = '  +0.57%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.96'
$ws.Range("E23").Value = '  +12.57%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.999'
$ws.Range("E24").Value = '  -0.12%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '65.81'
$ws.Range("E25").Value = '  -2.38%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '627.29'
$ws.Range("E26").Value = '  +13.00%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000103'
$ws.Range("E27").Value = '  +9.02%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.51'
$ws.Range("E28").Value = '  -1.75%  '

# Row 29
$ws.Range("D29").Value = '2.589.66'
$ws.Range("E29").Value = '  +1.61%  '

# Row 30
$ws.Range("B30").Value = 'Fetch.AI'
$ws.Range("C30").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.50'
$ws.Range("E30").Value = '  +5.74%  '

# Row 31
$ws.Range("B31").Value = 'Binance-PegBSC-USD'
$ws.Range("C31").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.998'
$ws.Range("E31").Value = '  -0.11%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.24'
$ws.Range("E32").Value = '  -0.53%  '

# Row 33
$ws.Range("B33").Value = 'Kaspa'
$ws.Range("C33").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.142'
$ws.Range("E33").Value = '  -3.90%  '

# Row 34
$ws.Range("B34").Value = 'PancakeSwap'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.90'
$ws.Range("E34").Value = '  +1.61%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.12'
$ws.Range("E35").Value = '  +6.08%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.49'
$ws.Range("E36").Value = '  -2.36%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.999'
$ws.Range("E37").Value = '  -0.03%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.383'
$ws.Range("E38").Value = '  +0.18%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.43'
$ws.Range("E39").Value = '  -2.59%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.78'
$ws.Range("E40").Value = '  +0.38%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '147.30'
$ws.Range("E41").Value = '  -1.87%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.79'
$ws.Range("E42").Value = '  -0.89%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.61'
$ws.Range("E43").Value = '  +13.02%  '

# Row 44
$ws.Range("E44").Value = '  -0.04%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '148.09'
$ws.Range("E45").Value = '  -0.14%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.73'
$ws.Range("E46").Value = '  +1.39%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '20.87'
$ws.Range("E47").Value = '  +3.06%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0540'
$ws.Range("E48").Value = '  +1.06%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.606'
$ws.Range("E49").Value = '  +1.50%  '

# Row 50
$ws.Range("E50").Value = '  +1.31%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0922'
$ws.Range("E51").Value = '  -0.35%  '
